$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 108, shifting rows 108:112 down to 109:113
$ws.Rows.Item(108).Insert()

# Populate the new row 108 with data
$ws.Cells.Item(108, 1).Value = 9
$ws.Cells.Item(108, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(108, 3).Value = "Metropolitana"
$ws.Cells.Item(108, 4).Value = 45147
$ws.Cells.Item(108, 5).Value = 13
$ws.Cells.Item(108, 6).Value = 100112035
$ws.Cells.Item(108, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 52
$ws.Cells.Item(108, 11).Value = 17000
$ws.Cells.Item(108, 12).Value = 18000
$ws.Cells.Item(108, 13).Value = 17500
$ws.Cells.Item(108, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(108, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(108, 16).Value = 1167
$ws.Cells.Item(108, 17).Value = 15
$ws.Cells.Item(108, 18).Value = "Hortaliza"
